$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cell F1 with the same style as E1 (bold, centered, bordered header style)
$ws.Range('F1').Value = 'time_taken'
$ws.Range('E1').Copy()
$ws.Range('F1').PasteSpecial(-4122)

# Populate F2:F198 with the time_taken values (plain, unstyled cells)
$timeValues = @(
    '2021-10-05 13:38:55.798880',
    '2021-10-05 13:38:55.798889',
    '2021-10-05 13:38:55.798892',
    '2021-10-05 13:38:55.798894',
    '2021-10-05 13:38:55.798896',
    '2021-10-05 13:38:55.798898',
    '2021-10-05 13:38:55.798900',
    '2021-10-05 13:38:55.798902',
    '2021-10-05 13:38:55.798904',
    '2021-10-05 13:38:55.798906',
    '2021-10-05 13:38:55.798908',
    '2021-10-05 13:38:55.798910',
    '2021-10-05 13:38:55.798912',
    '2021-10-05 13:38:55.798914',
    '2021-10-05 13:38:55.798916',
    '2021-10-05 13:38:55.798917',
    '2021-10-05 13:38:55.798920',
    '2021-10-05 13:38:55.798921',
    '2021-10-05 13:38:55.798923',
    '2021-10-05 13:38:55.798925',
    '2021-10-05 13:38:55.798927',
    '2021-10-05 13:38:55.798929',
    '2021-10-05 13:38:55.798931',
    '2021-10-05 13:38:55.798933',
    '2021-10-05 13:38:55.798935',
    '2021-10-05 13:38:55.798937',
    '2021-10-05 13:38:55.798939',
    '2021-10-05 13:38:55.798941',
    '2021-10-05 13:38:55.798942',
    '2021-10-05 13:38:55.798944',
    '2021-10-05 13:38:55.798946',
    '2021-10-05 13:38:55.798948',
    '2021-10-05 13:38:55.798950',
    '2021-10-05 13:38:55.798952',
    '2021-10-05 13:38:55.798954',
    '2021-10-05 13:38:55.798956',
    '2021-10-05 13:38:55.798957',
    '2021-10-05 13:38:55.798959',
    '2021-10-05 13:38:55.798961',
    '2021-10-05 13:38:55.798963',
    '2021-10-05 13:38:55.798965',
    '2021-10-05 13:38:55.798967',
    '2021-10-05 13:38:55.798969',
    '2021-10-05 13:38:55.798971',
    '2021-10-05 13:38:55.798973',
    '2021-10-05 13:38:55.798975',
    '2021-10-05 13:38:55.798977',
    '2021-10-05 13:38:55.798979',
    '2021-10-05 13:38:55.798980',
    '2021-10-05 13:38:55.798982',
    '2021-10-05 13:38:55.798984',
    '2021-10-05 13:38:55.798986',
    '2021-10-05 13:38:55.798988',
    '2021-10-05 13:38:55.798990',
    '2021-10-05 13:38:55.798992',
    '2021-10-05 13:38:55.798994',
    '2021-10-05 13:38:55.798996',
    '2021-10-05 13:38:55.798998',
    '2021-10-05 13:38:55.799000',
    '2021-10-05 13:38:55.799002',
    '2021-10-05 13:38:55.799004',
    '2021-10-05 13:38:55.799006',
    '2021-10-05 13:38:55.799007',
    '2021-10-05 13:38:55.799009',
    '2021-10-05 13:38:55.799012',
    '2021-10-05 13:38:55.799014',
    '2021-10-05 13:38:55.799016',
    '2021-10-05 13:38:55.799018',
    '2021-10-05 13:38:55.799020',
    '2021-10-05 13:38:55.799022',
    '2021-10-05 13:38:55.799023',
    '2021-10-05 13:38:55.799026',
    '2021-10-05 13:38:55.799028',
    '2021-10-05 13:38:55.799030',
    '2021-10-05 13:38:55.799031',
    '2021-10-05 13:38:55.799033',
    '2021-10-05 13:38:55.799037',
    '2021-10-05 13:38:55.799039',
    '2021-10-05 13:38:55.799041',
    '2021-10-05 13:38:55.799043',
    '2021-10-05 13:38:55.799045',
    '2021-10-05 13:38:55.799047',
    '2021-10-05 13:38:55.799049',
    '2021-10-05 13:38:55.799051',
    '2021-10-05 13:38:55.799052',
    '2021-10-05 13:38:55.799054',
    '2021-10-05 13:38:55.799056',
    '2021-10-05 13:38:55.799058',
    '2021-10-05 13:38:55.799060',
    '2021-10-05 13:38:55.799062',
    '2021-10-05 13:38:55.799064',
    '2021-10-05 13:38:55.799066',
    '2021-10-05 13:38:55.799068',
    '2021-10-05 13:38:55.799070',
    '2021-10-05 13:38:55.799072',
    '2021-10-05 13:38:55.799074',
    '2021-10-05 13:38:55.799076',
    '2021-10-05 13:38:55.799078',
    '2021-10-05 13:38:55.799080',
    '2021-10-05 13:38:55.799082',
    '2021-10-05 13:38:55.799084',
    '2021-10-05 13:38:55.799085',
    '2021-10-05 13:38:55.799088',
    '2021-10-05 13:38:55.799090',
    '2021-10-05 13:38:55.799092',
    '2021-10-05 13:38:55.799094',
    '2021-10-05 13:38:55.799096',
    '2021-10-05 13:38:55.799098',
    '2021-10-05 13:38:55.799101',
    '2021-10-05 13:38:55.799103',
    '2021-10-05 13:38:55.799105',
    '2021-10-05 13:38:55.799107',
    '2021-10-05 13:38:55.799109',
    '2021-10-05 13:38:55.799111',
    '2021-10-05 13:38:55.799113',
    '2021-10-05 13:38:55.799115',
    '2021-10-05 13:38:55.799116',
    '2021-10-05 13:38:55.799118',
    '2021-10-05 13:38:55.799120',
    '2021-10-05 13:38:55.799122',
    '2021-10-05 13:38:55.799124',
    '2021-10-05 13:38:55.799126',
    '2021-10-05 13:38:55.799127',
    '2021-10-05 13:38:55.799129',
    '2021-10-05 13:38:55.799131',
    '2021-10-05 13:38:55.799133',
    '2021-10-05 13:38:55.799135',
    '2021-10-05 13:38:55.799137',
    '2021-10-05 13:38:55.799140',
    '2021-10-05 13:38:55.799143',
    '2021-10-05 13:38:55.799157',
    '2021-10-05 13:38:55.799160',
    '2021-10-05 13:38:55.799163',
    '2021-10-05 13:38:55.799166',
    '2021-10-05 13:38:55.799169',
    '2021-10-05 13:38:55.799172',
    '2021-10-05 13:38:55.799175',
    '2021-10-05 13:38:55.799178',
    '2021-10-05 13:38:55.799180',
    '2021-10-05 13:38:55.799183',
    '2021-10-05 13:38:55.799186',
    '2021-10-05 13:38:55.799189',
    '2021-10-05 13:38:55.799192',
    '2021-10-05 13:38:55.799196',
    '2021-10-05 13:38:55.799198',
    '2021-10-05 13:38:55.799200',
    '2021-10-05 13:38:55.799202',
    '2021-10-05 13:38:55.799204',
    '2021-10-05 13:38:55.799206',
    '2021-10-05 13:38:55.799208',
    '2021-10-05 13:38:55.799210',
    '2021-10-05 13:38:55.799212',
    '2021-10-05 13:38:55.799214',
    '2021-10-05 13:38:55.799217',
    '2021-10-05 13:38:55.799220',
    '2021-10-05 13:38:55.799224',
    '2021-10-05 13:38:55.799227',
    '2021-10-05 13:38:55.799230',
    '2021-10-05 13:38:55.799232',
    '2021-10-05 13:38:55.799234',
    '2021-10-05 13:38:55.799236',
    '2021-10-05 13:38:55.799238',
    '2021-10-05 13:38:55.799240',
    '2021-10-05 13:38:55.799242',
    '2021-10-05 13:38:55.799244',
    '2021-10-05 13:38:55.799247',
    '2021-10-05 13:38:55.799250',
    '2021-10-05 13:38:55.799253',
    '2021-10-05 13:38:55.799256',
    '2021-10-05 13:38:55.799259',
    '2021-10-05 13:38:55.799262',
    '2021-10-05 13:38:55.799265',
    '2021-10-05 13:38:55.799271',
    '2021-10-05 13:38:55.799275',
    '2021-10-05 13:38:55.799279',
    '2021-10-05 13:38:55.799282',
    '2021-10-05 13:38:55.799285',
    '2021-10-05 13:38:55.799287',
    '2021-10-05 13:38:55.799289',
    '2021-10-05 13:38:55.799291',
    '2021-10-05 13:38:55.799294',
    '2021-10-05 13:38:55.799296',
    '2021-10-05 13:38:55.799298',
    '2021-10-05 13:38:55.799300',
    '2021-10-05 13:38:55.799302',
    '2021-10-05 13:38:55.799305',
    '2021-10-05 13:38:55.799307',
    '2021-10-05 13:38:55.799308',
    '2021-10-05 13:38:55.799310',
    '2021-10-05 13:38:55.799312',
    '2021-10-05 13:38:55.799314',
    '2021-10-05 13:38:55.799316',
    '2021-10-05 13:38:55.799318',
    '2021-10-05 13:38:55.799321',
    '2021-10-05 13:38:55.799323',
    '2021-10-05 13:38:55.799325',
    '2021-10-05 13:38:55.799327'
)

for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeValues[$i]
}
